function Set-CellXml($Row, $Col, $Xml) {
    $doc = $word.ActiveDocument
    $tbl = $doc.Tables.Item(3)
    $cell = $tbl.Cell($Row, $Col)
    $cell.Range.InsertXML($Xml)
    $doc2 = $word.ActiveDocument
    $cell2 = $doc2.Tables.Item(3).Cell($Row, $Col)
    $firstPara = $cell2.Range.Paragraphs.Item(1)
    $firstPara.Range.Delete()
}

$xml_row4_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Minta kiválasztása, csapattagok között a feladatok megbeszélése, véglegesítése</w:t></w:r></w:p>
'@

Set-CellXml 4 2 $xml_row4_col2

$xml_row9_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Elemzőket futtatni és ezek eredményt átnézni, javítani, majd erről dokumentációt írni</w:t></w:r></w:p>
'@

Set-CellXml 9 2 $xml_row9_col2

$xml_row14_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Tool</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-ok beüzemelése</w:t></w:r></w:p>
'@

Set-CellXml 14 2 $xml_row14_col2

$xml_row15_col2_merge = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Funkcionális tesztek megírásának elkezdése, mind emellett a monitorozást kezelni</w:t></w:r></w:p>
'@

Set-CellXml 15 2 $xml_row15_col2_merge

$xml_row19_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Tool</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-ok output-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jainak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> összehasonlítani, majd a hasonlóságokról dokumentációt írni</w:t></w:r></w:p>
'@

Set-CellXml 19 2 $xml_row19_col2

$xml_row24_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>H</w:t></w:r><w:r><w:t>eti tesztek hibamenedzselése</w:t></w:r></w:p>
'@

Set-CellXml 24 2 $xml_row24_col2

$xml_row29_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Utolsó tesztek hibamenedzselése</w:t></w:r></w:p>
'@

Set-CellXml 29 2 $xml_row29_col2

$xml_row34_col1_lrpb = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Török Zoltán</w:t></w:r></w:p>
'@

Set-CellXml 34 1 $xml_row34_col1_lrpb

$xml_row34_col2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">„Bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hunting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, dokumentáció, tesztek utolsó ellenőrzése</w:t></w:r><w:r><w:t xml:space="preserve">, esetleges </w:t></w:r><w:r><w:t>le</w:t></w:r><w:r><w:t>maradt rész</w:t></w:r><w:r><w:t>feladatba</w:t></w:r><w:r><w:t xml:space="preserve"> besegíteni</w:t></w:r></w:p>
'@

Set-CellXml 34 2 $xml_row34_col2

$xml_row35_col2_merge = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Esetlegesen, ha kell megtalált hibák javítása</w:t></w:r></w:p>
'@

Set-CellXml 35 2 $xml_row35_col2_merge

$xml_row37_col1_no_lrpb = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Egri Bálint</w:t></w:r></w:p>
'@

Set-CellXml 37 1 $xml_row37_col1_no_lrpb
